# Fruta / hortaliza, semanal
# A new weekly record was inserted as row 7 (pushing the previous rows 7-73
# down to rows 8-74), growing the sheet's used range from A1:R73 to A1:R74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7; everything below shifts down one row.
$ws.Rows("7:7").Insert()

# Populate the new row 7 with the new weekly record.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Macroferia Regional de Talca"
$ws.Range("C7").Value = "Maule"
$ws.Range("D7").Value = 44532
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 100112022
$ws.Range("G7").Value = "Arveja Verde"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 13000
$ws.Range("N7").Value = "$/saco 25 kilos"
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 520
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
